# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" sheets, matching the new scrape results.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Map of row -> new value for sheet "展览"
$sheet1Updates = @{
    2  = 2715
    13 = 9208
    17 = 260
    23 = 1000
    24 = 2096
    25 = 2190
    27 = 1891
    28 = 1930
    30 = 1547
    31 = 279
    32 = 160
    38 = 491
    39 = 12
    40 = 61
    41 = 629
    42 = 39
    43 = 1397
    44 = 304
    46 = 180
    47 = 646
    49 = 298
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Map of row -> new value for sheet "全部类型"
$sheet4Updates = @{
    2  = 2715
    10 = 9208
    16 = 260
    21 = 1000
    22 = 2190
    23 = 1891
    25 = 1547
    26 = 279
    27 = 160
    33 = 491
    37 = 12
    38 = 61
    39 = 629
    41 = 39
    42 = 1397
    44 = 304
    46 = 180
    47 = 646
    48 = 298
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}

$wb.Save()
